# =====================================================================
# Adds "Power Supply" and "Gain" worksheets with transformer / LM3886
# power-supply / gain calculations, plus a couple of new SPL rows and
# two new speaker-link hyperlinks on the existing "SPL" sheet.
# =====================================================================

$wb = $excel.ActiveWorkbook
$wsSPL = $wb.Worksheets.Item("SPL")

$xlCenter = -4108

# ---------------------------------------------------------------
# 1. Add the two new worksheets in the right order: SPL, Power Supply, Gain
# ---------------------------------------------------------------
$wsPower = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsPower.Name = "Power Supply"

$wsGain = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsGain.Name = "Gain"

# =====================================================================
# 2. "Power Supply" sheet
# =====================================================================

# --- Section 1: Transformer Specifications ---
$wsPower.Range("B1").Value = "Transformer Specifications"
$wsPower.Range("B1").HorizontalAlignment = $xlCenter
$wsPower.Range("B1").Font.Bold = $true
$wsPower.Range("C1").HorizontalAlignment = $xlCenter
$wsPower.Range("C1").Font.Bold = $true
$wsPower.Range("B1:C1").Merge()

$wsPower.Range("B2").Value = "Voltage Regulation Error"
$wsPower.Range("C2").Value = 0.07

$wsPower.Range("B4").Value = "Input (VAC)"
$wsPower.Range("C4").Value = "Output (VAC)"
$wsPower.Range("D4").Value = "Inductance Ratio"
$wsPower.Range("E4").Value = "Rectified DC Voltage"

$wsPower.Range("B5").Value = 115
$wsPower.Range("C5").Value = 24
$wsPower.Range("D5").Formula = "=POWER((B5/C5), 2)"
$wsPower.Range("E5").Formula = "=C5*SQRT(2)"

$wsPower.Range("B6").Value = 120
$wsPower.Range("C6").Formula = "=B6/SQRT(D6)"
$wsPower.Range("D6").Value = 22.96
$wsPower.Range("E6").Formula = "=C6*SQRT(2)"

# --- Section 2: LM3886 Power Supply Requirements ---
$wsPower.Range("B14").Value = "LM3886 Power Supply Requirements"
$wsPower.Range("B14").HorizontalAlignment = $xlCenter
$wsPower.Range("B14").Font.Bold = $true
$wsPower.Range("C14").HorizontalAlignment = $xlCenter
$wsPower.Range("C14").Font.Bold = $true
$wsPower.Range("B14:C14").Merge()

$wsPower.Range("C16").Value = "Load Impedance (" + [char]0x03A9 + ")"
$c16chars = $wsPower.Range("C16").Characters(17, 2)
$c16chars.Font.Name = "Calibri"

$wsPower.Range("D16").Value = "Average Output Power (W)"
$wsPower.Range("E16").Value = "Peak Output Voltage"
$wsPower.Range("F16").Value = "Peak Output Current"

$wsPower.Range("G16").Value = "Maximum Supply Voltage (" + [char]0x00B1 + "V)"
$g16chars = $wsPower.Range("G16").Characters(25, 3)
$g16chars.Font.Name = "Calibri"

$wsPower.Range("H16").Value = "Minimum Gain"

$wsPower.Range("J16").Value = "Transformer Voltage Regulation"
$wsPower.Range("K16").Value = 0.07

$wsPower.Range("B17").Value = "Tweeter"
$wsPower.Range("C17").Value = 4
$wsPower.Range("D17").Value = 50
$wsPower.Range("E17").Formula = "=SQRT(2*C17*D17)"
$wsPower.Range("F17").Formula = "=SQRT(2*D17/C17)"
$wsPower.Range("G17").Formula = "=(E17+K`$18)*(1+K`$16)*(1+K`$17)"
$wsPower.Range("H17").Formula = "=SQRT(D17*C17)"

$wsPower.Range("J17").Value = "Mains Voltage Variation"
$wsPower.Range("K17").Value = 0.1

$wsPower.Range("B18").Value = "Woofer"
$wsPower.Range("C18").Value = 4
$wsPower.Range("D18").Value = 60
$wsPower.Range("E18").Formula = "=SQRT(2*C18*D18)"
$wsPower.Range("F18").Formula = "=SQRT(2*D18/C18)"
$wsPower.Range("G18").Formula = "=(E18+K`$18)*(1+K`$16)*(1+K`$17)"
$wsPower.Range("H18").Formula = "=SQRT(D18*C18)"

$wsPower.Range("J18").Value = "Drop-out voltage of LM3886"
$wsPower.Range("K18").Value = 4

$wsPower.Range("J19").Value = "Input Voltage Level"
$wsPower.Range("K19").Value = 1

$wsPower.Columns("B:H").AutoFit() | Out-Null
$wsPower.Columns("J:K").AutoFit() | Out-Null

$wsPower.PageSetup.Orientation = 1

$wsPower.Range("D12").Select()

# =====================================================================
# 3. "Gain" sheet
# =====================================================================
$wsGain.Range("B1").Value = "Tweeters"
$wsGain.Range("B1").HorizontalAlignment = $xlCenter
$wsGain.Range("B1").Font.Bold = $true
$wsGain.Range("C1").HorizontalAlignment = $xlCenter
$wsGain.Range("C1").Font.Bold = $true
$wsGain.Range("B1:C1").Merge()

$wsGain.Range("E8").Select()

# =====================================================================
# 4. "SPL" sheet additions
# =====================================================================
$wsSPL.Range("A6").Value = "X"
$wsSPL.Range("A7").Value = "X"

$wsSPL.Range("D7").Value = 1
$wsSPL.Range("E7").Value = 106.1

$wsSPL.Range("B8").Value = 93.4
$wsSPL.Range("C8").Value = 50
$wsSPL.Range("D8").Value = 1
$wsSPL.Range("E8").Value = 110.4

$wsSPL.Range("F8").Value = "https://www.parts-express.com/tc-6024-6-1-2-treated-paper-cone-woofer-with-foam-surround-4-ohm--299-2196"
$wsSPL.Hyperlinks.Add($wsSPL.Range("F8"), "https://www.parts-express.com/tc-6024-6-1-2-treated-paper-cone-woofer-with-foam-surround-4-ohm--299-2196") | Out-Null
$wsSPL.Range("F8").Style = $wsSPL.Range("F6").Style

$wsSPL.Range("F9").Value = "https://www.parts-express.com/goldwood-gw-8024-8-butyl-surround-woofer-4-ohm--290-356"
$wsSPL.Hyperlinks.Add($wsSPL.Range("F9"), "https://www.parts-express.com/goldwood-gw-8024-8-butyl-surround-woofer-4-ohm--290-356") | Out-Null
$wsSPL.Range("F9").Style = $wsSPL.Range("F6").Style

$wsSPL.Range("C12").Select()

# =====================================================================
# 5. Final workbook state: Gain tab active, as in the authored workbook
# =====================================================================
$wsGain.Activate()
$wsGain.Range("E8").Select()

Write-Host "Edit complete"
